$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cell H1 (bold, centered, thin border).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new I and J columns of data (values are identical per row).
$data = @{
    2  = 8
    3  = 8
    4  = 8
    5  = 7
    6  = 8
    7  = 8
    8  = 8
    9  = 6
    10 = 7
    11 = 8
    12 = 8
    13 = 8
    14 = 7
    15 = 8
    16 = 9
    17 = 8
    18 = 8
    19 = 8
    20 = 8
    21 = 8
    22 = 8
    23 = 8
    24 = 8
    25 = 8
    26 = 8
    27 = 9
    28 = 8
    29 = 3
    30 = 5
    31 = 5
    32 = 3
}

foreach ($row in $data.Keys) {
    $value = $data[$row]
    $ws.Cells.Item($row, 9).Value = $value   # column I
    $ws.Cells.Item($row, 10).Value = $value  # column J
}
